$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "\href{https://rr.peercommunityin.org/}{Peer Community In Registered Reports}"

# Row 8 holds the "Journals Include" section header (columns A-D) together
# with the first journal bullet in column E. Insert a blank row below the
# header row so the existing bullet list (rows 9 onward) shifts down by
# one row, then move the previous first bullet into the new row and place
# the new entry at the top of the list (still on row 8).
$ws.Rows.Item(9).Insert()
$ws.Range("E9").Value = $ws.Range("E8").Value()
$ws.Range("E8").Value = $newText

# Replace the short PCI name used in the summary row with the full name.
$ws.Range("C2").Value = $newText

# Update the row 2 height to accommodate the longer wrapped text.
$ws.Rows.Item(2).RowHeight = 45

# Update the active selection to reflect the editor's final cursor position.
$ws.Range("C2").Select()
